$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E4").Value = 69.13
$ws1.Range("D7").Value = 468.29
$ws1.Range("H7").Value = 1073.7
$ws1.Range("M16").Value = 6512.36
$ws1.Range("M20").Value = 307.93
$ws1.Range("I23").Value = 26.1
$ws1.Range("M23").Value = 383.47
$ws1.Range("D25").Value = 2289.6
$ws1.Range("C33").Value = 1555.2
$ws1.Range("D33").Value = 2799.36
$ws1.Range("L33").Value = 855.36
$ws1.Range("O37").Value = 4136.14
$ws1.Range("E51").Value = 142.11
$ws1.Range("H51").Value = 180
$ws1.Range("D55").Value = 570.24
$ws1.Range("E55").Value = 55.65
$ws1.Range("M55").Value = 777.39

# Row 58 "x de 56" summary counters (text cells)
$ws1.Range("C58").Value = "3 de 56"
$ws1.Range("D58").Value = "10 de 56"
$ws1.Range("E58").Value = "4 de 56"
$ws1.Range("H58").Value = "3 de 56"
$ws1.Range("I58").Value = "4 de 56"
$ws1.Range("L58").Value = "5 de 56"
$ws1.Range("M58").Value = "14 de 56"
$ws1.Range("O58").Value = "5 de 56"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 373.31
$ws2.Range("F7").Value = 1541.99
$ws2.Range("F16").Value = 8065.41
$ws2.Range("F20").Value = 307.93
$ws2.Range("F23").Value = 409.57
$ws2.Range("F25").Value = 6346.89
$ws2.Range("F33").Value = 5209.92
$ws2.Range("F37").Value = 15847.07
$ws2.Range("F51").Value = 322.11
$ws2.Range("F55").Value = 2063.26
$ws2.Range("F58").Value = 62018.24

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D2").Value = 3529.27
$ws3.Range("E2").Value = 6441.07304517915
$ws3.Range("F2").Value = 0.3539767873590337
$ws3.Range("D3").Value = 16832.44
$ws3.Range("E3").Value = 10624.5676
$ws3.Range("F3").Value = 0.6130471406505347
$ws3.Range("D4").Value = 802.67
$ws3.Range("E4").Value = 200.33
$ws3.Range("F4").Value = 0.8002691924227318
$ws3.Range("D7").Value = 1521
$ws3.Range("E7").Value = 879
$ws3.Range("F7").Value = 0.63375
$ws3.Range("D8").Value = 756.66
$ws3.Range("E8").Value = 243.34
$ws3.Range("F8").Value = 0.75666
$ws3.Range("D15").Value = 6163.29
$ws3.Range("E15").Value = 7336.71
$ws3.Range("F15").Value = 0.45654
$ws3.Range("D16").Value = 24512.01
$ws3.Range("E16").Value = 31547.69
$ws3.Range("F16").Value = 0.4372483263378149
$ws3.Range("D18").Value = 6070.05
$ws3.Range("E18").Value = -2870.05
$ws3.Range("F18").Value = 1.896890625
$ws3.Range("D19").Value = 62018.24
$ws3.Range("E19").Value = 55421.45064517915
$ws3.Range("F19").Value = 0.5280858597233186
